# USER MASTER.xlsx - password comparison refactor
#
# The "Password" column of Table1 (sheet "User") used to store a dummy
# numeric placeholder (123) for every user. As part of the rework of the
# password comparison logic, the placeholder is replaced with a literal
# text value "abc123" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")
$lo = $ws.ListObjects.Item("Table1")

$passwordColumn = $lo.ListColumns.Item("Password")
$dataRange = $passwordColumn.DataBodyRange

for ($i = 1; $i -le $dataRange.Rows.Count; $i++) {
    $dataRange.Cells.Item($i, 1).Value = "abc123"
}

# Leave the selection where the edit session ended up.
$ws.Range("D5").Select()
